# Sonoverse workbook update: insert a new "Aerobilia (Pneumobilia)" row
# into the Gallbladder and biliary tract section (new row 5), shifting
# all subsequent rows down by one, and refresh the hyperlinks/selection
# accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank row at row 5 (pushes former rows 5-25 to 6-26).
$ws.Rows.Item(5).Insert()

# 2. Populate the new row 5 with the new term.
$ws.Range("A5").Value2 = "Gallbladder and biliary tract"
$ws.Range("B5").Value2 = "Aerobilia (Pneumobilia)"
$ws.Range("C5").Value2 = "Clip 1 B-mode"
$ws.Range("D5").Value2 = "https://youtu.be/K2Wbg7BgXy4 "

# Copy the hyperlink cell style (index 1, "Collegamento ipertestuale") from
# a neighboring hyperlink cell so the new D5 looks the same as the rest.
$ws.Range("D5").Style = $ws.Range("D6").Style

# 3. Rebuild the hyperlinks collection: the row insert does not remap the
#    worksheet's existing hyperlink ranges automatically, so drop them all
#    and re-add them pointing at their now-shifted rows, in the same
#    relative order as before, then append the brand-new one last.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("D3"), "https://youtu.be/zxTC0YBY2RY ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D20"), "https://youtu.be/xBfd04F4Ni8 ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D10"), "https://youtu.be/91M82AIMyu0 ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D24"), "https://youtu.be/qushjTAy6XQ ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D22"), "https://youtu.be/pc-vbxSRTbs ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D16"), "https://youtu.be/DjI1kEnzfSQ ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D21"), "https://youtu.be/JvwODCASLYQ ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D17"), "https://youtu.be/U3ydTsRwxok ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D12"), "https://youtu.be/15o_Km86IzM ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D25"), "https://youtu.be/_FckFwJwynI ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D23"), "https://youtu.be/Axbee4vjNtU") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D14"), "https://youtu.be/RhSUFLTmTl4") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D6"), "https://youtu.be/2kRZcpi70Aw ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D26"), "https://youtu.be/z_oaRVxRz5s ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D5"), "https://youtu.be/K2Wbg7BgXy4 ") | Out-Null

# 4. Update the sort state range to cover the newly-added row.
$ws.Sort.SetRange($ws.Range("A2:C24"))
$ws.Sort.SortFields.Item(1).SetRange($ws.Range("A2:A24"))

# 5. Move the active selection to where the author left off (D7).
$ws.Range("D7").Select()
